# unit tests + colC
# Replicates: new template banner text, merged A1:I1 header, collapsed
# column A:B width, cleared A2, refreshed selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the banner text in A1, drop the old "email address" cell (A2) ---
$ws.Range("A1").Value = "This is the template sheet for the package LStest – if you do not want to see this…"
$ws.Range("A2").Clear()

# --- Row 1 grows a touch taller and becomes an explicit (custom) height ---
$ws.Rows.Item(1).RowHeight = 27.6

# --- Columns A and B collapse into a single ~11.5-wide column band ---
$ws.Range("A:B").ColumnWidth = 10.69

# --- Merge the header across A1:I1 (this also extends the s="2" style
#     of A1 across B1:I1, matching the blank styled cells in the target) ---
$ws.Range("A1:I1").Merge()

# --- Refresh the saved selection/active cell ---
$ws.Range("C8").Select() | Out-Null

Write-Host "edit applied"
